$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B7").Value = "SingleUseId8"
$ws.Range("C7").Value = "Iceland_45"
$ws.Range("D7").Value = "Left"
$ws.Range("E7").Value = "LTR"
$ws.Range("F7").Value = "PWM & ACEL"
